$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear cell C2 entirely (remove the cell, not just its content)
$ws.Range("C2").ClearContents()

# Update the active selection to C2
$ws.Range("C2").Select()
